$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(6, 6).Value = 2
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(15, 6).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(23, 6).Value = 1
$ws.Cells.Item(26, 6).Value = 3
$ws.Cells.Item(40, 6).Value = 2
$ws.Cells.Item(42, 6).Value = 1
$ws.Cells.Item(52, 6).Value = 1
$ws.Cells.Item(53, 6).Value = 2
$ws.Cells.Item(55, 6).Value = 1
$ws.Cells.Item(56, 6).Value = 1
$ws.Cells.Item(58, 6).Value = 1
$ws.Cells.Item(64, 6).Value = 2
$ws.Cells.Item(65, 6).Value = 3
$ws.Cells.Item(68, 6).Value = 1
$ws.Cells.Item(71, 6).Value = 1
$ws.Cells.Item(75, 6).Value = 1
$ws.Cells.Item(78, 6).Value = 1
$ws.Cells.Item(93, 6).Value = 1
$ws.Cells.Item(95, 6).Value = 1
$ws.Cells.Item(96, 6).Value = 2
$ws.Cells.Item(97, 6).Value = 2
$ws.Cells.Item(98, 6).Value = 2
$ws.Cells.Item(99, 6).Value = 1
$ws.Cells.Item(102, 6).Value = 1
$ws.Cells.Item(103, 6).Value = 2
$ws.Cells.Item(104, 6).Value = 1
$ws.Cells.Item(105, 6).Value = 1
$ws.Cells.Item(106, 6).Value = 1
$ws.Cells.Item(111, 6).Value = 1
$ws.Cells.Item(118, 6).Value = 1
$ws.Cells.Item(120, 6).Value = 1
$ws.Cells.Item(133, 6).Value = 1
$ws.Cells.Item(138, 6).Value = 1
$ws.Cells.Item(139, 6).Value = 1
$ws.Cells.Item(167, 6).Value = 1
$ws.Cells.Item(173, 6).Value = 1
$ws.Cells.Item(184, 6).Value = 1
$ws.Cells.Item(213, 6).Value = 2
$ws.Cells.Item(214, 6).Value = 1
$ws.Cells.Item(215, 6).Value = 1
$ws.Cells.Item(216, 6).Value = 1
$ws.Cells.Item(220, 6).Value = 1
$ws.Cells.Item(221, 6).Value = 1
$ws.Cells.Item(223, 6).Value = 2
$ws.Cells.Item(231, 6).Value = 1
$ws.Cells.Item(246, 6).Value = 1
$ws.Cells.Item(248, 6).Value = 1
$ws.Cells.Item(252, 6).Value = 2
$ws.Cells.Item(253, 6).Value = 1
$ws.Cells.Item(258, 6).Value = 1
$ws.Cells.Item(259, 6).Value = 1
$ws.Cells.Item(264, 6).Value = 1
$ws.Cells.Item(284, 6).Value = 1
$ws.Cells.Item(291, 6).Value = 1
$ws.Cells.Item(298, 6).Value = 1
$ws.Cells.Item(302, 6).Value = 1
$ws.Cells.Item(303, 6).Value = 1
$ws.Cells.Item(305, 6).Value = 1
$ws.Cells.Item(306, 6).Value = 1
$ws.Cells.Item(307, 6).Value = 2
$ws.Cells.Item(311, 6).Value = 1
